$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.096.65"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "3.476.06"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.88%  "
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "4.067.80"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "3.476.06"
$ws.Range("E15").Value = "  -0.33%  "
$ws.Range("D16").Value = "64.079.70"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "384.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.567"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "3.616.30"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("E31").Value = "  +2.87%  "
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "3.503.44"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0778"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.76%  "
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "2.330.90"
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("E51").Value = "  -2.80%  "
